# Auto-generated edit script: updates Leve profit-calculation sheets
# with refreshed market-price data (scheduled runner update).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 446.625
$ws.Range("I55").Value = 516.8333
$ws.Range("J55").Value = 236
$ws.Range("K55").Value = 516.8333
$ws.Range("L55").Value = 236
$ws.Range("M55").Value = -302.8333
$ws.Range("N55").Value = -664
$ws.Range("H98").Value = 1834.3334
$ws.Range("I98").Value = 1834.3334
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1834.3334
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -336.3334
$ws.Range("H112").Value = 21164946
$ws.Range("J112").Value = 24845738
$ws.Range("L112").Value = 74537214
$ws.Range("N112").Value = -74539430
$ws.Range("H122").Value = 1834.3334
$ws.Range("I122").Value = 1834.3334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5503.0002
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -3053.0002
$ws.Range("H124").Value = 25750
$ws.Range("J124").Value = 25750
$ws.Range("L124").Value = 25750
$ws.Range("N124").Value = -35570
$ws.Range("H137").Value = 2203.8215
$ws.Range("I137").Value = 1495.8182
$ws.Range("J137").Value = 4799.8335
$ws.Range("K137").Value = 4487.4546
$ws.Range("L137").Value = 14399.5005
$ws.Range("M137").Value = -1937.4546
$ws.Range("N137").Value = -19499.5005
# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 10445.5
$ws.Range("I43").Value = 6421
$ws.Range("J43").Value = 12457.75
$ws.Range("K43").Value = 6421
$ws.Range("L43").Value = 12457.75
$ws.Range("M43").Value = -6108
$ws.Range("N43").Value = -13083.75
$ws.Range("H45").Value = 18292.334
$ws.Range("I45").Value = 18292.334
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 18292.334
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -17915.334
$ws.Range("H74").Value = 1427.5682
$ws.Range("I74").Value = 1099.4688
$ws.Range("J74").Value = 2302.5
$ws.Range("K74").Value = 1099.4688
$ws.Range("L74").Value = 2302.5
$ws.Range("M74").Value = -225.4688000000001
$ws.Range("N74").Value = -4050.5
$ws.Range("H77").Value = 1427.5682
$ws.Range("I77").Value = 1099.4688
$ws.Range("J77").Value = 2302.5
$ws.Range("K77").Value = 5497.344000000001
$ws.Range("L77").Value = 11512.5
$ws.Range("M77").Value = -1129.344000000001
$ws.Range("N77").Value = -20248.5
$ws.Range("H82").Value = 42000
$ws.Range("J82").Value = 42000
$ws.Range("L82").Value = 42000
$ws.Range("N82").Value = -42722
$ws.Range("H85").Value = 42000
$ws.Range("J85").Value = 42000
$ws.Range("L85").Value = 42000
$ws.Range("N85").Value = -44496
$ws.Range("H122").Value = 1224209.4
$ws.Range("I122").Value = 1352778.8
$ws.Range("K122").Value = 4058336.4
$ws.Range("M122").Value = -4055886.4
$ws.Range("H132").Value = 2463.7966
$ws.Range("I132").Value = 1294.5555
$ws.Range("J132").Value = 4293.913
$ws.Range("K132").Value = 3883.6665
$ws.Range("L132").Value = 12881.739
$ws.Range("M132").Value = -1353.6665
$ws.Range("N132").Value = -17941.739
# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1652.075
$ws.Range("I132").Value = 1160.8077
$ws.Range("J132").Value = 2564.4285
$ws.Range("K132").Value = 3482.4231
$ws.Range("L132").Value = 7693.2855
$ws.Range("M132").Value = -952.4231
$ws.Range("N132").Value = -12753.2855
# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 360
$ws.Range("I18").Value = 360
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1080
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -911
$ws.Range("H41").Value = 2900
$ws.Range("J41").Value = 2900
$ws.Range("L41").Value = 8700
$ws.Range("N41").Value = -9376
$ws.Range("H116").Value = 1212.1111
$ws.Range("I116").Value = 844.1429000000001
$ws.Range("K116").Value = 2532.4287
$ws.Range("M116").Value = 909.5712999999996
$ws.Range("H124").Value = 7666.3335
$ws.Range("J124").Value = 10000
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820
$ws.Range("H130").Value = 4042.9453
$ws.Range("I130").Value = 2286.6667
$ws.Range("J130").Value = 4118.2144
$ws.Range("K130").Value = 6860.000100000001
$ws.Range("L130").Value = 12354.6432
$ws.Range("M130").Value = -1840.000100000001
$ws.Range("N130").Value = -22394.6432
$ws.Range("H131").Value = 1667684.8
$ws.Range("I131").Value = 11111900
$ws.Range("J131").Value = 1058.5098
$ws.Range("K131").Value = 33335700
$ws.Range("L131").Value = 3175.5294
$ws.Range("M131").Value = -33330660
$ws.Range("N131").Value = -13255.5294
$ws.Range("H133").Value = 47799.582
$ws.Range("J133").Value = 6437.5
$ws.Range("L133").Value = 19312.5
$ws.Range("N133").Value = -29432.5
# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2550.6667
$ws.Range("I80").Value = 2486.353
$ws.Range("J80").Value = 2660
$ws.Range("K80").Value = 2486.353
$ws.Range("L80").Value = 2660
$ws.Range("M80").Value = -1488.353
$ws.Range("N80").Value = -4656
$ws.Range("H83").Value = 2550.6667
$ws.Range("I83").Value = 2486.353
$ws.Range("J83").Value = 2660
$ws.Range("K83").Value = 12431.765
$ws.Range("L83").Value = 13300
$ws.Range("M83").Value = -7439.764999999999
$ws.Range("N83").Value = -23284
$ws.Range("H122").Value = 23151914
$ws.Range("I122").Value = 3970753.8
$ws.Range("J122").Value = 50005540
$ws.Range("K122").Value = 11912261.4
$ws.Range("L122").Value = 150016620
$ws.Range("M122").Value = -11909811.4
$ws.Range("N122").Value = -150021520
$ws.Range("H132").Value = 20427.527
$ws.Range("I132").Value = 30607.828
$ws.Range("K132").Value = 91823.484
$ws.Range("M132").Value = -89293.484
# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 55557868
$ws.Range("I40").Value = 83335880
$ws.Range("J40").Value = 1834.1666
$ws.Range("K40").Value = 83335880
$ws.Range("L40").Value = 1834.1666
$ws.Range("M40").Value = -83335744
$ws.Range("N40").Value = -2106.1666
$ws.Range("H122").Value = 2911607.8
$ws.Range("I122").Value = 3250001
$ws.Range("J122").Value = 1670833.4
$ws.Range("K122").Value = 9750003
$ws.Range("L122").Value = 5012500.199999999
$ws.Range("M122").Value = -9747553
$ws.Range("N122").Value = -5017400.199999999
# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2030
$ws.Range("I122").Value = 2287.5
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 6862.5
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -4412.5
$ws.Range("N122").Value = -7900
